$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q3" sheet right before the existing "2022-Q2"
#    sheet. We clone the "2022-Q2" sheet itself so the new sheet starts
#    out with identical column layout / cell styling (bold+bordered
#    header row, bold+bordered index column) and then overwrite its
#    cell values with the 2022-Q3 fund-holdings data.
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$existingQ2.Copy($existingQ2)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# The source sheet had 18 data rows (rows 2-18); the new sheet only
# needs 2 data rows, so drop the extra rows the copy brought along.
$newSheet.Range("4:18").Delete()

# Columns B (基金代码) and D:G (基金规模/股票总仓位/仓位占比/持有市值) hold
# numeric-looking values that must stay TEXT (leading zeros / exact
# decimal strings), so force a text number format before writing them.
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "011410"
$newSheet.Range("C2").Value = "中信建投量化进取6个月持有期混合A"
$newSheet.Range("D2").Value = "5.36"
$newSheet.Range("E2").Value = "90.16"
$newSheet.Range("F2").Value = "0.62"
$newSheet.Range("G2").Value = "0.0332"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "011411"
$newSheet.Range("C3").Value = "中信建投量化进取6个月持有期混合C"
$newSheet.Range("D3").Value = "1.63"
$newSheet.Range("E3").Value = "90.16"
$newSheet.Range("F3").Value = "0.62"
$newSheet.Range("G3").Value = "0.0101"
$newSheet.Range("H3").Value = 10

# Drop the forced text format back to the sheet's normal style so the
# cells don't carry a stray numFmt (matches the plain un-styled data
# cells used throughout the rest of the workbook).
$newSheet.Range("B2:B3").Style = "Normal"
$newSheet.Range("D2:G3").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: add a 2022-Q3 row at the top of
#    the data (row 2) and shift the rest down by one, appending the
#    2020-Q4 row that falls off the end at row 9.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$rows = @(
    @("2022-Q3", 2, 0.04),
    @("2022-Q2", 17, 2.07),
    @("2022-Q1", 20, 4.41),
    @("2021-Q4", 19, 13.16),
    @("2021-Q3", 26, 11.79),
    @("2021-Q2", 9, 2.35),
    @("2021-Q1", 12, 2.87),
    @("2020-Q4", 3, 0.82)
)

# Row 9 doesn't exist yet, so clone the style of row 8's index cell
# (bold + border) onto it before the value-writing loop below.
$totalSheet.Range("A8").Copy($totalSheet.Range("A9"))

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $rows[$i][0]
    $totalSheet.Cells.Item($r, 3).Value = $rows[$i][1]
    $totalSheet.Cells.Item($r, 4).Value = $rows[$i][2]
}

# Restore the original active sheet/selection (sheet copy/insert above
# shifts focus onto the newly created sheet).
$totalSheet.Activate() | Out-Null
$totalSheet.Range("A1").Select() | Out-Null
